# Insert a new row at position 402, shifting existing rows 402-478 down to 403-479,
# then populate the new row 402 with the new data record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(402).Insert()

$ws.Cells.Item(402, 1).Value = 11
$ws.Cells.Item(402, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(402, 3).Value = "Bíobío"
$ws.Cells.Item(402, 4).Value = 45015
$ws.Cells.Item(402, 5).Value = 8
$ws.Cells.Item(402, 6).Value = 100114014
$ws.Cells.Item(402, 7).Value = "Betarraga"
$ws.Cells.Item(402, 8).Value = "Sin especificar"
$ws.Cells.Item(402, 9).Value = "Primera"
$ws.Cells.Item(402, 10).Value = 900
$ws.Cells.Item(402, 11).Value = 600
$ws.Cells.Item(402, 12).Value = 650
$ws.Cells.Item(402, 13).Value = 622
$ws.Cells.Item(402, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(402, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(402, 16).Value = 124
$ws.Cells.Item(402, 17).Value = 5
$ws.Cells.Item(402, 18).Value = "Hortaliza"

$ws.Cells.Item(402, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
